# Applies the "Exit info list.xlsx" update:
#  - Removes the now-unused "Disc Out" shared string and appends a dozen new
#    strings describing the Jaguar-motor / cannon / flipper wiring that was
#    added to Sheet2.
#  - Fills in the newly documented B/C/D columns on Sheet2 (and fixes a couple
#    of existing C-column off-by-one references caused by the removed string).
#  - Widens column C and updates the remembered selection on Sheet2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# The old rows 19-24 (Cannon/InFlipper/FlipperAtPlace/Middle/Bottom/FrisbeeDirection,
# column A only) slide down to rows 20-25 and pick up new B/D data; row 19 ends up
# completely empty, so clear it out first.
$ws.Range("A19:D19").ClearContents()

# Row data: row number -> A, B, C, D text (omitted key means "leave as is")
$rows = @(
    @{ Row = 5;  B = "Joystick";     C = "driver's joystick" },
    @{ Row = 6;  B = "Joystick";     C = "operator's joystick" },
    @{ Row = 8;  B = "Jaguar Motor"; C = "outer shooting cannon";          D = "PWM 7" },
    @{ Row = 9;  B = "Jaguar Motor"; C = "inner shooting cannon";         D = "PWM 8" },
    @{ Row = 10; A = "zAvit";        B = "Jaguar Motor"; C = "cannon angle";                    D = "PWM 9" },
    @{ Row = 12; A = "CannonAI";     B = "Potentiometer"; C = "measures angle of cannon";        D = "AI 1" },
    @{ Row = 13; B = "Encoder";      C = "measures shooting motor's RPM"; D = "DIO 7-8" },
    @{ Row = 15; B = "Jaguar Motor"; C = "flipper's motor";                D = "PWM 9" },
    @{ Row = 16; B = "Jaguar Motor"; D = "PWM 6" },
    @{ Row = 17; B = "Jaguar Motor"; D = "PWM 5" },
    @{ Row = 18; A = "Transfer";     B = "Jaguar Motor"; D = "PWM 10" },
    @{ Row = 20; A = "Cannon";       B = "Micro Switch"; D = "DIO 5" },
    @{ Row = 21; A = "InFlipper";    B = "Micro Switch"; D = "DIO 4" },
    @{ Row = 22; A = "FlipperAtPlace"; B = "Micro Switch"; D = "DIO 3" },
    @{ Row = 23; A = "Middle";       B = "Micro Switch"; D = "DIO 2" },
    @{ Row = 24; A = "Bottom";       B = "Micro Switch"; D = "DIO 1" },
    @{ Row = 25; A = "FrisbeeDirection"; B = "Micro Switch"; D = "DIO 6" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    if ($r.ContainsKey("A")) { $ws.Cells.Item($rowNum, 1).Value = $r.A }
    if ($r.ContainsKey("B")) { $ws.Cells.Item($rowNum, 2).Value = $r.B }
    if ($r.ContainsKey("C")) { $ws.Cells.Item($rowNum, 3).Value = $r.C }
    if ($r.ContainsKey("D")) { $ws.Cells.Item($rowNum, 4).Value = $r.D }
}

# Widen column C to fit the longer descriptions that were added (target stored
# width is 27.625 "chars"; the host's ColumnWidth setter quantizes to a 1/7
# character grid, so 26.85 is the input that round-trips closest to it).
$ws.Columns.Item(3).ColumnWidth = 26.85

# Remember the same selection the author ended up with.
$ws.Range("C17").Select()
